# Weekly data refresh: insert the latest week's Kiwi price record for
# "Vega Monumental Concepción" as a new row right after the existing
# row 234 block, pushing all following rows down by one (A1:T332 -> A1:T333).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 235; everything currently at/after 235
# (rows 235..332) shifts down to 236..333.
$ws.Rows.Item(235).Insert()

# Populate the newly inserted row 235 with this week's record.
$ws.Range("A235").Value = 11
$ws.Range("B235").Value = "Vega Monumental Concepción"
$ws.Range("C235").Value = "Bíobío"
$ws.Range("D235").Value = 45146
$ws.Range("E235").Value = 8
$ws.Range("F235").Value = "Fruta"
$ws.Range("G235").Value = 100101
$ws.Range("H235").Value = "Berries"
$ws.Range("I235").Value = 100101007
$ws.Range("J235").Value = "Kiwi"
$ws.Range("K235").Value = "Hayward"
$ws.Range("L235").Value = "Primera"
$ws.Range("M235").Value = 220
$ws.Range("N235").Value = 13000
$ws.Range("O235").Value = 14000
$ws.Range("P235").Value = 13455
$ws.Range("Q235").Value = "$/bandeja 18 kilos"
$ws.Range("R235").Value = "Región de O'Higgins"
$ws.Range("S235").Value = 748
$ws.Range("T235").Value = 18
